# Generate Report for Handoff
# b.md has been handed off again (new source/target files + new handoff datetime).
# Update the "Overview" sheet status, and the zh-cn / de-de detail sheets
# (Status, Latest Handoff File, Latest Handoff Datetime for row 3 = b.md).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Overview sheet: status for b.md (row 3) changes from
# "Handed back: in sync with en-US" to "Ready for handoff" for both locales.
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

# ---------------------------------------------------------------------------
# zh-cn sheet: row 3 (b.md) gets a new handoff.
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = "Ready for handoff"
$zhcn.Range("C3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("D3").Value = "2016-02-23 07:30:30"

# Rebuild the zh-cn hyperlinks collection so the new display text on C3 is
# picked up (existing external hyperlinks are read-only for in-place edits).
$zhcnLinks = @(
    @{ Ref = "A2"; Address = "https://github.com/OpenLocalizationTest/oltest/blob/aca4bc4c16767daee601f6e0722e72996287eda2/e2e/a.md"; Display = "a.md" },
    @{ Ref = "C2"; Address = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e539883cfd49214c05963e347f9c26d185e9b2ca/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/terryjin/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"; Display = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf" },
    @{ Ref = "E2"; Address = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/c519afee8f4106816ba158eaa4fe05a8f7e4a798/e2e/a.md"; Display = "a.md" },
    @{ Ref = "F2"; Address = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/0e1bdfff3bfeed95b6df070cc940196b5dd2a7ed/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/terryjin/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"; Display = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf" },
    @{ Ref = "A3"; Address = "https://github.com/OpenLocalizationTest/oltest/blob/aca4bc4c16767daee601f6e0722e72996287eda2/e2e/b.md"; Display = "b.md" },
    @{ Ref = "C3"; Address = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e539883cfd49214c05963e347f9c26d185e9b2ca/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/terryjin/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"; Display = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf" },
    @{ Ref = "E3"; Address = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/c519afee8f4106816ba158eaa4fe05a8f7e4a798/e2e/a.md"; Display = "a.md" },
    @{ Ref = "F3"; Address = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/0e1bdfff3bfeed95b6df070cc940196b5dd2a7ed/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/terryjin/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"; Display = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf" },
    @{ Ref = "A4"; Address = "https://github.com/OpenLocalizationTest/oltest/blob/aca4bc4c16767daee601f6e0722e72996287eda2/.localization-config"; Display = ".localization-config" }
)

$zhcn.Hyperlinks.Delete()
foreach ($link in $zhcnLinks) {
    $zhcn.Hyperlinks.Add($zhcn.Range($link.Ref), $link.Address, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $link.Display)
}

# ---------------------------------------------------------------------------
# de-de sheet: row 3 (b.md) gets a new handoff.
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = "Ready for handoff"
$dede.Range("C3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("D3").Value = "2016-02-23 07:30:47"

$dedeLinks = @(
    @{ Ref = "A2"; Address = "https://github.com/OpenLocalizationTest/oltest/blob/aca4bc4c16767daee601f6e0722e72996287eda2/e2e/a.md"; Display = "a.md" },
    @{ Ref = "C2"; Address = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/baeb1a34094d07e71a8ac46d838a16fd45085b00/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/terryjin/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"; Display = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf" },
    @{ Ref = "E2"; Address = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/1d23e265ec11abcc170313efb6e3777c8b836b99/e2e/a.md"; Display = "a.md" },
    @{ Ref = "F2"; Address = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/3495a52f94e6294cc8ee5323a3c7307e1cfc2781/ol-handback/OpenLocalizationTestOrg/oltest.de-de/terryjin/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"; Display = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf" },
    @{ Ref = "A3"; Address = "https://github.com/OpenLocalizationTest/oltest/blob/aca4bc4c16767daee601f6e0722e72996287eda2/e2e/b.md"; Display = "b.md" },
    @{ Ref = "C3"; Address = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/baeb1a34094d07e71a8ac46d838a16fd45085b00/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/terryjin/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"; Display = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf" },
    @{ Ref = "E3"; Address = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/1d23e265ec11abcc170313efb6e3777c8b836b99/e2e/a.md"; Display = "a.md" },
    @{ Ref = "F3"; Address = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/3495a52f94e6294cc8ee5323a3c7307e1cfc2781/ol-handback/OpenLocalizationTestOrg/oltest.de-de/terryjin/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"; Display = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf" },
    @{ Ref = "A4"; Address = "https://github.com/OpenLocalizationTest/oltest/blob/aca4bc4c16767daee601f6e0722e72996287eda2/.localization-config"; Display = ".localization-config" }
)

$dede.Hyperlinks.Delete()
foreach ($link in $dedeLinks) {
    $dede.Hyperlinks.Add($dede.Range($link.Ref), $link.Address, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $link.Display)
}

Write-Host "Handoff report generated for b.md"
